$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the comment on E1 (was attached to the "seat_min" header before the
# column was renamed/repurposed).
if ($ws.Range("E1").Comment) {
    $ws.Range("E1").Comment.Delete()
}

# Rename header E1 from "seat_min" to "seat_last"
$ws.Range("E1").Value = "seat_last"

# Row2 (was seat_min=1, seat_max=14) becomes seat_last = "r2s7"
$ws.Range("E2").Value = "r2s7"

# Row3 (was seat_min=1, seat_max=6) becomes seat_last = "r2s3"
$ws.Range("E3").Value = "r2s3"

# Select column F, then delete it (old "seat_max" column), matching the
# end-user gesture of selecting the column before removing it.
[void]$ws.Columns("F").Select()
$ws.Columns("F").Delete()
